$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.213.83"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "1.795.25"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.41"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4486"
$ws.Range("E7").Value = "  +15.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3716"
$ws.Range("E8").Value = "  +10.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.14"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07560"
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.35"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.290"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.460"
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("D16").Value = "1.792.25"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06740"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.11"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +4.50%  "
$ws.Range("E22").Value = "  +3.53%  "
$ws.Range("D23").Value = "28.199.34"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.78"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.48"
$ws.Range("E26").Value = "  +4.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.97"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.361"
$ws.Range("D29").Value = "1.998.41"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.99"
$ws.Range("E30").Value = "  +4.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.240"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.040"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09390"
$ws.Range("E33").Value = "  +7.88%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.793"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2365"
$ws.Range("E35").Value = "  +13.86%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02333"
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.208"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6552"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.330"
$ws.Range("E41").Value = "  +6.61%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.479"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.25"
$ws.Range("E44").Value = "  +5.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.832"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6073"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.80"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.024"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07124"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.159"
$ws.Range("E51").Value = "  +1.89%  "
